$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 74
$ws.Range("F3").Value = 160
$ws.Range("F5").Value = 45
$ws.Range("F8").Value = 1609
$ws.Range("F9").Value = 7388
$ws.Range("F11").Value = 7564
$ws.Range("F12").Value = 15
$ws.Range("F13").Value = 32
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 6041
$ws.Range("F16").Value = 3230
$ws.Range("F17").Value = 3597
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 10
$ws.Range("F21").Value = 24
$ws.Range("F24").Value = 275
$ws.Range("F26").Value = 2083
$ws.Range("F30").Value = 252
$ws.Range("F31").Value = 1040
$ws.Range("F32").Value = 59
$ws.Range("F33").Value = 6
$ws.Range("F34").Value = 2583
$ws.Range("F35").Value = 1418
$ws.Range("F36").Value = 4
$ws.Range("F37").Value = 2
$ws.Range("F38").Value = 8
$ws.Range("F39").Value = 3174
$ws.Range("F40").Value = 140
$ws.Range("F41").Value = 232
$ws.Range("F44").Value = 471
$ws.Range("F45").Value = 1231
$ws.Range("F46").Value = 223
$ws.Range("F47").Value = 517
$ws.Range("F48").Value = 579

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 52
$ws.Range("F9").Value = 392

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 74
$ws.Range("F5").Value = 160
$ws.Range("F7").Value = 52
$ws.Range("F8").Value = 45
$ws.Range("F11").Value = 1609
$ws.Range("F14").Value = 7388
$ws.Range("F16").Value = 7564
$ws.Range("F17").Value = 15
$ws.Range("F18").Value = 6041
$ws.Range("F19").Value = 3230
$ws.Range("F20").Value = 3597
$ws.Range("F21").Value = 10
$ws.Range("F22").Value = 10
$ws.Range("F23").Value = 24
$ws.Range("F26").Value = 275
$ws.Range("F30").Value = 2083
$ws.Range("F36").Value = 1040
$ws.Range("F37").Value = 59
$ws.Range("F38").Value = 2583
$ws.Range("F39").Value = 1418
$ws.Range("F41").Value = 3174
$ws.Range("F42").Value = 140
$ws.Range("F45").Value = 471
$ws.Range("F46").Value = 1231
$ws.Range("F47").Value = 223
$ws.Range("F48").Value = 517
$ws.Range("F49").Value = 579
